$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowValues = @(45763,485,481,518,538,550.25,540,504,488,487,709,719,481,481,557,559,486,625,486,490,540,450,450,477,476,565,495,490,484,580,471,497,547,495,499,530,526.5,548,545,500,814,598,613,612,569,593,640)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(27, $i + 1).Value = $rowValues[$i]
}

$ws.Range("A27").NumberFormat = "YYYY-MM-DD"
